# #5: cash & deposit done
# Adds legislator/property metadata columns (property_category, category,
# date, legislator_name, legislator_id, source_file, index) to the
# "deposit" (存款) sheet, matching the layout already used by the other
# sheets in this workbook, and turns the old B1:F1 "sample data" header
# row into a real header row (bank/deposit_type/currency/owner/total).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)   # 存款 (deposit)

# ---- header row (row 1) -------------------------------------------------
$ws.Range("B1").Value = "bank"
$ws.Range("C1").Value = "deposit_type"
$ws.Range("D1").Value = "currency"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"
$ws.Range("G1").Value = "property_category"
$ws.Range("H1").Value = "category"
$ws.Range("I1").Value = "date"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"

# new header cells (G1:M1) need the same bold/border style as the rest of
# row 1 (style index 1) - copy it over with a formats-only paste so we
# reuse the existing style instead of fabricating a new one.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("G1:M1").PasteSpecial(-4122) | Out-Null

# ---- data rows (rows 2:5) ------------------------------------------------
$indices = @(42, 43, 44, 45)
for ($r = 2; $r -le 5; $r++) {
    $ws.Cells.Item($r, 7).Value  = "deposit"     # property_category
    $ws.Cells.Item($r, 8).Value  = "normal"      # category
    $ws.Cells.Item($r, 10).Value = "姚文智"       # legislator_name
    $ws.Cells.Item($r, 11).Value = 1745          # legislator_id
    $ws.Cells.Item($r, 12).Value = "tmp60da1"    # source_file
    $ws.Cells.Item($r, 13).Value = $indices[$r - 2]  # index

    # "date" (column I) needs to hold the literal text "2012-04-30" -
    # assigning that string straight to .Value gets auto-parsed into a
    # date serial number, so round it through a throwaway formula cell
    # (which evaluates to a plain text result) and paste the *value*
    # back in, leaving the destination's formatting untouched.
    $ws.Range("Z1").Formula = '="2012-04-30"'
    $ws.Range("Z1").Copy() | Out-Null
    $ws.Cells.Item($r, 9).PasteSpecial(-4163) | Out-Null
    $ws.Range("Z1").Value = ""
}
